$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 92; this shifts existing rows 92-112 down to 93-113,
# matching the diff (dimension A1:R112 -> A1:R113, with a brand-new weekly price
# record inserted at row 92).
$ws.Rows.Item(92).Insert()

# Populate the new row 92 with the new weekly record.
$ws.Range("A92").Value = 7
$ws.Range("B92").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C92").Value = "Ñuble"
$ws.Range("D92").Value = 44543
$ws.Range("E92").Value = 16
$ws.Range("F92").Value = 100112024
$ws.Range("G92").Value = "Choclo"
$ws.Range("H92").Value = "Choclero"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 16000
$ws.Range("K92").Value = 400
$ws.Range("L92").Value = 450
$ws.Range("M92").Value = 425
$ws.Range("N92").Value = "$/unidad"
$ws.Range("O92").Value = "Región de O'Higgins"
$ws.Range("P92").Value = 425
$ws.Range("Q92").Value = 1
$ws.Range("R92").Value = "Hortaliza"
